$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 277, shifting the existing rows 277:373 down to 278:374.
$ws.Rows.Item(277).Insert()

# Populate the newly inserted row 277 with the new weekly record.
$ws.Range("A277").Value2 = 7
$ws.Range("B277").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C277").Value2 = "Ñuble"
$ws.Range("D277").Value2 = 44468
$ws.Range("E277").Value2 = 16
$ws.Range("F277").Value2 = 100112033
$ws.Range("G277").Value2 = "Lechuga"
$ws.Range("H277").Value2 = "Escarola"
$ws.Range("I277").Value2 = "Primera"
$ws.Range("J277").Value2 = 300
$ws.Range("K277").Value2 = 9000
$ws.Range("L277").Value2 = 10000
$ws.Range("M277").Value2 = 9500
$ws.Range("N277").Value2 = "`$/caja 15 unidades"
$ws.Range("O277").Value2 = "Provincia del Elquí"
$ws.Range("P277").Value2 = 633
$ws.Range("Q277").Value2 = 15
$ws.Range("R277").Value2 = "Hortaliza"

# Make sure the style (date number format) on the new date cell matches the rest of column D.
$ws.Range("D277").NumberFormat = $ws.Range("D278").NumberFormat
